$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$rowData = @{
  27 = @(3.1800595238095242, 2.9940476190476191, 3.1934523809523809, 3.0610119047619051, 3.041666666666667, 2.8258928571428572, 2.816964285714286, 2.8973214285714279, 2.796130952380953, 2.8184523809523809, 2.59375, 2.8556547619047619, 2.8214285714285721, 2.8363095238095242, 2.808383233532934, 3.15625)
  28 = @(3.860215053763441, 3.471774193548387, 3.84005376344086, 3.763440860215054, 3.635752688172043, 3.450268817204301, 3.307795698924731, 3.588709677419355, 3.19758064516129, 3.400537634408602, 2.9946236559139781, 3.556451612903226, 3.344086021505376, 3.385752688172043, 3.3766937669376689, 3.637096774193548)
  29 = @(3.6541666666666668, 3.4874999999999998, 3.7361111111111112, 3.6694444444444438, 3.6430555555555562, 3.458333333333333, 3.2527777777777782, 3.5277777777777781, 3.4361111111111109, 3.4777777777777779, 3.119444444444444, 3.568055555555556, 3.3708333333333331, 3.4486111111111111, 3.347826086956522, 3.5222222222222221)
  30 = @(3.370967741935484, 3.205645161290323, 3.5053763440860219, 3.341397849462366, 3.342741935483871, 3.268817204301075, 3.1411290322580649, 3.259408602150538, 3.091397849462366, 3.143817204301075, 2.758064516129032, 3.225806451612903, 3.143817204301075, 3.205645161290323, 3.1972789115646258, 3.341397849462366)
  31 = @(2.2124999999999999, 2.1138888888888889, 2.3611111111111112, 2.255555555555556, 2.208333333333333, 2.2777777777777781, 2.1833333333333331, 2.0972222222222219, 2.0555555555555549, 2.0819444444444439, 1.8638888888888889, 2.072222222222222, 2.0541666666666671, 2.0180555555555562, 1.9282639885222379, 2.3736111111111109)
  32 = @(2.692204301075269, 2.416666666666667, 2.774193548387097, 2.665322580645161, 2.635752688172043, 2.625, 2.475806451612903, 2.497311827956989, 2.491935483870968, 2.4946236559139781, 2.189516129032258, 2.6088709677419351, 2.514784946236559, 2.684139784946237, 2.3799185888738128, 2.84005376344086)
  33 = @(2.770161290322581, 2.4516129032258061, 2.870967741935484, 2.5981182795698921, 2.6733870967741939, 2.55241935483871, 2.46505376344086, 2.501344086021505, 2.342741935483871, 2.638440860215054, 2.25, 2.567204301075269, 2.495967741935484, 2.579301075268817, 2.458108108108108, 2.940860215053763)
  34 = @(5.4722222222222223, 4.8250000000000002, 5.1833333333333336, 4.8361111111111112, 5.1472222222222221, 4.697222222222222, 4.7874999999999996, 4.7263888888888888, 4.7833333333333332, 5.1958333333333337, 4.3944444444444448, 4.7680555555555557, 4.7374999999999998, 4.7166666666666668, 4.6215083798882679, 5.0236111111111112)
  35 = @(4.198924731182796, 4.096774193548387, 4.309139784946237, 4.014784946236559, 4.17741935483871, 3.854838709677419, 3.873655913978495, 3.888440860215054, 4.344086021505376, 4.372311827956989, 3.708333333333333, 3.958333333333333, 3.970430107526882, 4.049731182795699, 4.0197740112994351, 4.115591397849462)
  36 = @(3.2875000000000001, 3.286111111111111, 3.5597222222222218, 3.1638888888888892, 3.3236111111111111, 2.9152777777777779, 2.9805555555555561, 3.161111111111111, 3.1597222222222219, 3.1652777777777779, 2.770833333333333, 3.1069444444444438, 3.05, 3.0874999999999999, 3.0154061624649859, 3.3944444444444439)
  37 = @(3.342741935483871, 3.381720430107527, 3.512096774193548, 3.192204301075269, 3.30510752688172, 3.111559139784946, 3.036290322580645, 3.158602150537634, 3.543010752688172, 3.452956989247312, 2.834677419354839, 3.03494623655914, 3.057795698924731, 3.125, 3.008163265306123, 3.666666666666667)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 2).Value = $vals[$i]
    }
}

# Update the view selection to match the saved state (scrolled to row 19, active cell D40)
[void]$ws.Range("D40").Select()

Write-Host "AQHI monthly average data completed"
